# Apply updated percentage values to the team matrix on Sheet1
# (changes to team matrices from games pulled march 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.2229508196721312
    $ws.Range("C2").Value = 0.5114754098360655
    $ws.Range("J2").Value = 0.01311475409836066
    $ws.Range("P2").Value = 0.1442622950819672
    $ws.Range("S2").Value = 0.1081967213114754
    $ws.Range("B3").Value = 0.006289308176100629
    $ws.Range("C3").Value = 0.02515723270440252
    $ws.Range("J3").Value = 0.03773584905660377
    $ws.Range("P3").Value = 0.7169811320754716
    $ws.Range("S3").Value = 0.2138364779874214
    $ws.Range("J4").Value = 0.02173913043478261
    $ws.Range("P4").Value = 0.7608695652173914
    $ws.Range("S4").Value = 0.2173913043478261
    $ws.Range("B6").Value = 0.05381165919282511
    $ws.Range("D6").Value = 0.02690582959641256
    $ws.Range("F6").Value = 0.03139013452914798
    $ws.Range("J6").Value = 0.2690582959641256
    $ws.Range("O6").Value = 0.008968609865470852
    $ws.Range("Q6").Value = 0.1479820627802691
    $ws.Range("R6").Value = 0.07174887892376682
    $ws.Range("S6").Value = 0.3901345291479821
    $ws.Range("B7").Value = 0.1005291005291005
    $ws.Range("D7").Value = 0.005291005291005291
    $ws.Range("F7").Value = 0.03174603174603174
    $ws.Range("J7").Value = 0.1322751322751323
    $ws.Range("O7").Value = 0.005291005291005291
    $ws.Range("Q7").Value = 0.1957671957671958
    $ws.Range("R7").Value = 0.08994708994708994
    $ws.Range("S7").Value = 0.4391534391534391
    $ws.Range("B8").Value = 0.09018036072144289
    $ws.Range("D8").Value = 0.01402805611222445
    $ws.Range("F8").Value = 0.05210420841683366
    $ws.Range("J8").Value = 0.1082164328657315
    $ws.Range("O8").Value = 0.01202404809619238
    $ws.Range("Q8").Value = 0.186372745490982
    $ws.Range("R8").Value = 0.09819639278557114
    $ws.Range("S8").Value = 0.438877755511022
    $ws.Range("B9").Value = 0.09448818897637795
    $ws.Range("D9").Value = 0.01181102362204724
    $ws.Range("F9").Value = 0.06692913385826772
    $ws.Range("J9").Value = 0.1062992125984252
    $ws.Range("O9").Value = 0.007874015748031496
    $ws.Range("Q9").Value = 0.1811023622047244
    $ws.Range("R9").Value = 0.1062992125984252
    $ws.Range("S9").Value = 0.4251968503937008
    $ws.Range("B10").Value = 0.09777777777777778
    $ws.Range("D10").Value = 0.02222222222222222
    $ws.Range("F10").Value = 0.06814814814814815
    $ws.Range("J10").Value = 0.1259259259259259
    $ws.Range("O10").Value = 0.00962962962962963
    $ws.Range("Q10").Value = 0.2096296296296296
    $ws.Range("R10").Value = 0.09185185185185185
    $ws.Range("S10").Value = 0.3748148148148148
    $ws.Range("G11").Value = 0.1312056737588652
    $ws.Range("J11").Value = 0.07446808510638298
    $ws.Range("K11").Value = 0.1666666666666667
    $ws.Range("L11").Value = 0.6099290780141844
    $ws.Range("S11").Value = 0.01773049645390071
    $ws.Range("G12").Value = 0.8023255813953488
    $ws.Range("J12").Value = 0.1686046511627907
    $ws.Range("K12").Value = 0.005813953488372093
    $ws.Range("L12").Value = 0.01162790697674419
    $ws.Range("S12").Value = 0.01162790697674419
    $ws.Range("G13").Value = 0.5675675675675675
    $ws.Range("J13").Value = 0.3513513513513514
    $ws.Range("S13").Value = 0.08108108108108109
    $ws.Range("F15").Value = 0.025
    $ws.Range("H15").Value = 0.175
    $ws.Range("I15").Value = 0.1
    $ws.Range("J15").Value = 0.315
    $ws.Range("K15").Value = 0.1
    $ws.Range("M15").Value = 0.02
    $ws.Range("O15").Value = 0.045
    $ws.Range("S15").Value = 0.22
    $ws.Range("F16").Value = 0.02105263157894737
    $ws.Range("H16").Value = 0.1947368421052632
    $ws.Range("I16").Value = 0.1
    $ws.Range("J16").Value = 0.3842105263157894
    $ws.Range("K16").Value = 0.09473684210526316
    $ws.Range("M16").Value = 0.01052631578947368
    $ws.Range("O16").Value = 0.06315789473684211
    $ws.Range("S16").Value = 0.131578947368421
    $ws.Range("F17").Value = 0.01848049281314168
    $ws.Range("H17").Value = 0.1601642710472279
    $ws.Range("I17").Value = 0.1006160164271047
    $ws.Range("J17").Value = 0.4517453798767967
    $ws.Range("K17").Value = 0.06981519507186858
    $ws.Range("M17").Value = 0.01642710472279261
    $ws.Range("N17").Value = 0.002053388090349076
    $ws.Range("O17").Value = 0.06570841889117043
    $ws.Range("S17").Value = 0.1149897330595483
    $ws.Range("F18").Value = 0.01282051282051282
    $ws.Range("H18").Value = 0.1965811965811966
    $ws.Range("I18").Value = 0.09401709401709402
    $ws.Range("J18").Value = 0.452991452991453
    $ws.Range("K18").Value = 0.1025641025641026
    $ws.Range("M18").Value = 0.0170940170940171
    $ws.Range("O18").Value = 0.02136752136752137
    $ws.Range("S18").Value = 0.1025641025641026
    $ws.Range("F19").Value = 0.01590744757772957
    $ws.Range("H19").Value = 0.2255965292841648
    $ws.Range("I19").Value = 0.1041214750542299
    $ws.Range("J19").Value = 0.3571945046999277
    $ws.Range("K19").Value = 0.09544468546637744
    $ws.Range("M19").Value = 0.01663051337671728
    $ws.Range("N19").Value = 0.0007230657989877079
    $ws.Range("O19").Value = 0.06507592190889371
    $ws.Range("S19").Value = 0.1193058568329718
